# Fixed missing department/region leaders
# - Insert a new "Regional Leader" row (2009, 15, 1, 16, 6%) right before the
#   existing "Attorney" row.
# - Append a new "Department Leader" row (2009, 1, 0, 1, 0%) at the end.
# - Re-state the termination-rate percentages without the decimal (100.0% ->
#   100%, 50.0% -> 50%, 23.5% -> 24%, 5.3% -> 5%).
# - Grow Table3 (and its AutoFilter) to cover the two extra rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Make room for the new "Regional Leader" row just above the "Attorney" row
# (current row 5) -- this pushes Attorney down to row 6 and keeps its cell
# types/values intact.
$ws.Rows.Item(5).Insert()

# --- New row 5: Regional Leader -----------------------------------------
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2009"
$ws.Range("A5").Style = "Normal"

$ws.Range("B5").Value = "Regional Leader"

$ws.Range("C5").Value = 15
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 16

$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "6%"
$ws.Range("F5").Style = "Normal"

# --- Row 2 (Product Manager): 100.0% -> 100% -----------------------------
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "100%"
$ws.Range("F2").Style = "Normal"

# --- Row 3 (Project Manager): 50.0% -> 50% -------------------------------
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "50%"
$ws.Range("F3").Style = "Normal"

# --- Row 4 (Paralegal): 23.5% -> 24% -------------------------------------
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "24%"
$ws.Range("F4").Style = "Normal"

# --- Row 6 (Attorney, was row 5): 5.3% -> 5% -----------------------------
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "5%"
$ws.Range("F6").Style = "Normal"

# --- New row 7: Department Leader ----------------------------------------
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2009"
$ws.Range("A7").Style = "Normal"

$ws.Range("B7").Value = "Department Leader"

$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1

$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "0%"
$ws.Range("F7").Style = "Normal"

# Grow the table (and its autofilter) so it covers rows 1-7.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F7"))
